$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting
# the existing "Late" / heading / "Outstanding" columns one to the right.
$ws.Columns("N:N").Insert()

# The newly inserted column should take the width of its left neighbour
# (column M, "In Advance") rather than keep any best-fit sizing.
$ws.Columns("N:N").ColumnWidth = 10.17

# Reflect the new used range and make this the active/selected sheet+cell,
# matching where the author left the selection after editing.
$ws.Activate()
$ws.Range("R7").Select()
